# Add new columns I ("I0") and J ("IF") to Sheet1, mirroring the
# existing header style (from H1) and filling in the computed values:
#   I = 1 for every data row
#   J = same value as column H for that row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine last used data row from column A (header in row 1, data starts row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 38 }

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (style) from H1 to I1:J1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats = -4122
$excel.CutCopyMode = $false

# Restore the values after pasting formats (PasteSpecial formats only
# shouldn't touch values, but keep this explicit for safety)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in data rows 2..lastRow
for ($r = 2; $r -le $lastRow; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2  # column H = 8 (numeric value)
    $ws.Cells.Item($r, 9).Value2 = 1       # column I = 9
    $ws.Cells.Item($r, 10).Value2 = $hVal  # column J = 10
}
